# Applies the commit "Minor revisions of templates" to
# poa_revocation_property.docx.
#
# The underlying XML diff is almost entirely Word re-flowing existing
# runs around newly inserted <w:proofErr/> spell/grammar-check markers:
# for every such hunk the concatenated <w:t> text is byte-for-byte
# identical before and after the change (those markers carry no visible
# content and are only ever produced as a side effect of the interactive
# spell/grammar checker, not through any documented, scriptable part of
# the Word object model). The one substantive, content-visible change in
# the diff is a new, empty "Default"-styled paragraph inserted right
# after the "Date: ___" paragraph, immediately before the paragraph that
# holds "{{ property_agent.name.full(...) }}".

$d = $word.ActiveDocument

# Locate the paragraph containing the property_agent name merge field;
# the new blank paragraph needs to be inserted immediately before it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*property_agent.name.full*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the property_agent.name.full paragraph"
}

# Insert a new, empty paragraph right before it. After this call,
# $target's Range refers to the newly created (now-preceding) empty
# paragraph rather than the original one.
$target.Range.InsertParagraphBefore()

# Give the new blank paragraph the "Default" style, matching the other
# blank spacer paragraphs used throughout this document.
$target.Style = "Default"
